$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '30.470.24'
$ws.Range('E2').Value = '  -1.40%  '

Set-TextValue 'D3' '2.097.28'
$ws.Range('E3').Value = '  -0.80%  '

$ws.Range('E4').Value = '  -0.22%  '

Set-TextValue 'D5' '330.82'
$ws.Range('E5').Value = '  -1.40%  '

Set-TextValue 'D6' '1.002'
$ws.Range('E6').Value = '  -0.08%  '

Set-TextValue 'D7' '0.5222'
$ws.Range('E7').Value = '  -1.71%  '

Set-TextValue 'D8' '0.4441'
$ws.Range('E8').Value = '  +1.60%  '

Set-TextValue 'D9' '54.00'
$ws.Range('E9').Value = '  +16.80%  '

Set-TextValue 'D10' '0.08930'
$ws.Range('E10').Value = '  -1.15%  '

Set-TextValue 'D11' '1.154'
$ws.Range('E11').Value = '  -1.96%  '

Set-TextValue 'D12' '24.43'
$ws.Range('E12').Value = '  -2.77%  '

Set-TextValue 'D13' '2.082.85'
$ws.Range('E13').Value = '  -1.46%  '

Set-TextValue 'D14' '6.697'
$ws.Range('E14').Value = '  -1.15%  '

Set-TextValue 'D15' '7.708'

Set-TextValue 'D16' '96.40'
$ws.Range('E16').Value = '  -1.27%  '

Set-TextValue 'D17' '1.004'
$ws.Range('E17').Value = '  -0.06%  '

$ws.Range('E18').Value = '  -1.02%  '

Set-TextValue 'D19' '0.06615'
$ws.Range('E19').Value = '  -0.86%  '

Set-TextValue 'D20' '19.15'
$ws.Range('E20').Value = '  -0.05%  '

$ws.Range('E21').Value = '  -0.22%  '

Set-TextValue 'D22' '6.277'
$ws.Range('E22').Value = '  -1.45%  '

Set-TextValue 'D23' '30.502.29'
$ws.Range('E23').Value = '  -1.54%  '

Set-TextValue 'D24' '12.29'
$ws.Range('E24').Value = '  +0.98%  '

$ws.Range('E25').Value = '  +1.93%  '

Set-TextValue 'D26' '2.338.24'
$ws.Range('E26').Value = '  -1.16%  '

$ws.Range('E27').Value = '  -2.30%  '

Set-TextValue 'D28' '2.570'
$ws.Range('E28').Value = '  +0.17%  '

Set-TextValue 'D29' '163.52'
$ws.Range('E29').Value = '  +0.12%  '

Set-TextValue 'D30' '132.07'
$ws.Range('E30').Value = '  -1.37%  '

Set-TextValue 'D31' '1.193'
$ws.Range('E31').Value = '  +1.88%  '

$ws.Range('E32').Value = '  -0.26%  '

Set-TextValue 'D33' '1.662'
$ws.Range('E33').Value = '  +9.07%  '

Set-TextValue 'D34' '6.168'
$ws.Range('E34').Value = '  -1.04%  '

Set-TextValue 'D35' '3.900'
$ws.Range('E35').Value = '  -2.80%  '

Set-TextValue 'D36' '10.21'
$ws.Range('E36').Value = '  +7.14%  '

Set-TextValue 'D37' '0.02567'
$ws.Range('E37').Value = '  -1.87%  '

Set-TextValue 'D38' '0.06816'
$ws.Range('E38').Value = '  +1.29%  '

Set-TextValue 'D39' '12.76'
$ws.Range('E39').Value = '  -1.47%  '

$ws.Range('E40').Value = '  -1.39%  '

Set-TextValue 'D41' '0.2263'
$ws.Range('E41').Value = '  -0.58%  '

Set-TextValue 'D42' '0.6914'
$ws.Range('E42').Value = '  +0.66%  '

Set-TextValue 'D43' '1.250'
$ws.Range('E43').Value = '  -0.10%  '

Set-TextValue 'D44' '1.001'
$ws.Range('E44').Value = '  -0.15%  '

Set-TextValue 'D45' '14.03'
$ws.Range('E45').Value = '  -0.71%  '

Set-TextValue 'D46' '0.6351'
$ws.Range('E46').Value = '  -1.57%  '

Set-TextValue 'D47' '2.250'
$ws.Range('E47').Value = '  +0.15%  '

Set-TextValue 'D48' '3.632'
$ws.Range('E48').Value = '  -1.43%  '

Set-TextValue 'D49' '1.244'
$ws.Range('E49').Value = '  +6.31%  '

$ws.Range('E50').Value = '  -2.44%  '

Set-TextValue 'D51' '81.91'
$ws.Range('E51').Value = '  -1.19%  '
